$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 154.28572
$ws.Range("I39").Value = 66.09999999999999
$ws.Range("J39").Value = 374.75
$ws.Range("K39").Value = 198.3
$ws.Range("L39").Value = 1124.25
$ws.Range("M39").Value = 97.70000000000002
$ws.Range("N39").Value = -1716.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1894.4445
$ws.Range("I86").Value = 1814.5
$ws.Range("J86").Value = 2174.25
$ws.Range("K86").Value = 1814.5
$ws.Range("L86").Value = 2174.25
$ws.Range("M86").Value = -691.5
$ws.Range("N86").Value = -4420.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1894.4445
$ws.Range("I89").Value = 1814.5
$ws.Range("J89").Value = 2174.25
$ws.Range("K89").Value = 9072.5
$ws.Range("L89").Value = 10871.25
$ws.Range("M89").Value = -3456.5
$ws.Range("N89").Value = -22103.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2333.3333
$ws.Range("I98").Value = 2333.3333
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2333.3333
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -835.3332999999998
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2628.8
$ws.Range("I125").Value = 732
$ws.Range("J125").Value = 3103
$ws.Range("K125").Value = 6588
$ws.Range("L125").Value = 27927
$ws.Range("M125").Value = -4128
$ws.Range("N125").Value = -32847

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3560.158
$ws.Range("I138").Value = 1654.8667
$ws.Range("J138").Value = 4028.672
$ws.Range("K138").Value = 4964.6001
$ws.Range("L138").Value = 12086.016
$ws.Range("M138").Value = 175.3999000000003
$ws.Range("N138").Value = -22366.016

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 1789.1111
$ws.Range("I39").Value = 1789.1111
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1789.1111
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1269.1111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1771.125
$ws.Range("I63").Value = 1433.8
$ws.Range("J63").Value = 2333.3333
$ws.Range("K63").Value = 1433.8
$ws.Range("L63").Value = 2333.3333
$ws.Range("M63").Value = -747.8
$ws.Range("N63").Value = -3705.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1771.125
$ws.Range("I66").Value = 1433.8
$ws.Range("J66").Value = 2333.3333
$ws.Range("K66").Value = 7169
$ws.Range("L66").Value = 11666.6665
$ws.Range("M66").Value = -3737
$ws.Range("N66").Value = -18530.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1291.2106
$ws.Range("I74").Value = 1348.3334
$ws.Range("J74").Value = 1193.2858
$ws.Range("K74").Value = 1348.3334
$ws.Range("L74").Value = 1193.2858
$ws.Range("M74").Value = -474.3334
$ws.Range("N74").Value = -2941.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1291.2106
$ws.Range("I77").Value = 1348.3334
$ws.Range("J77").Value = 1193.2858
$ws.Range("K77").Value = 6741.666999999999
$ws.Range("L77").Value = 5966.429
$ws.Range("M77").Value = -2373.666999999999
$ws.Range("N77").Value = -14702.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 966.6667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 966.6667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 966.6667
$ws.Range("N88").Value = -1778.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 966.6667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 966.6667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 966.6667
$ws.Range("N91").Value = -3774.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1731.1111
$ws.Range("I122").Value = 1722.5
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5167.5
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2717.5
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 32500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 32500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 32500
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -33332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2765.73
$ws.Range("I105").Value = 1841.2858
$ws.Range("J105").Value = 2835.3118
$ws.Range("K105").Value = 1841.2858
$ws.Range("L105").Value = 2835.3118
$ws.Range("M105").Value = -94.28580000000011
$ws.Range("N105").Value = -6329.311799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1680.4651
$ws.Range("I134").Value = 1618.3939
$ws.Range("J134").Value = 1885.3
$ws.Range("K134").Value = 4855.1817
$ws.Range("L134").Value = 5655.9
$ws.Range("M134").Value = -2320.1817
$ws.Range("N134").Value = -10725.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5001227
$ws.Range("I6").Value = 8750022
$ws.Range("J6").Value = 2833.3333
$ws.Range("K6").Value = 8750022
$ws.Range("L6").Value = 2833.3333
$ws.Range("M6").Value = -8749909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 500.875
$ws.Range("I10").Value = 500.875
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500.875
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -361.875
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1585187.4
$ws.Range("I19").Value = 2375093.5
$ws.Range("J19").Value = 5375
$ws.Range("K19").Value = 2375093.5
$ws.Range("L19").Value = 5375
$ws.Range("M19").Value = -2374923.5
$ws.Range("N19").Value = -5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 1585187.4
$ws.Range("I24").Value = 2375093.5
$ws.Range("J24").Value = 5375
$ws.Range("K24").Value = 2375093.5
$ws.Range("L24").Value = 5375
$ws.Range("M24").Value = -2374923.5
$ws.Range("N24").Value = -5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2310.5386
$ws.Range("I31").Value = 1949.579
$ws.Range("J31").Value = 3290.2856
$ws.Range("K31").Value = 1949.579
$ws.Range("L31").Value = 3290.2856
$ws.Range("M31").Value = -1654.579
$ws.Range("N31").Value = -3880.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2310.5386
$ws.Range("I34").Value = 1949.579
$ws.Range("J34").Value = 3290.2856
$ws.Range("K34").Value = 1949.579
$ws.Range("L34").Value = 3290.2856
$ws.Range("M34").Value = -1747.579
$ws.Range("N34").Value = -3694.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1659.6364
$ws.Range("I122").Value = 1341
$ws.Range("J122").Value = 1841.7142
$ws.Range("K122").Value = 4023
$ws.Range("L122").Value = 5525.142599999999
$ws.Range("M122").Value = -1573
$ws.Range("N122").Value = -10425.1426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2881.5518
$ws.Range("I132").Value = 1891.1111
$ws.Range("J132").Value = 4502.273
$ws.Range("K132").Value = 5673.3333
$ws.Range("L132").Value = 13506.819
$ws.Range("M132").Value = -3143.3333
$ws.Range("N132").Value = -18566.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2293.5483
$ws.Range("I134").Value = 1411.1666
$ws.Range("J134").Value = 5318.857
$ws.Range("K134").Value = 4233.4998
$ws.Range("L134").Value = 15956.571
$ws.Range("M134").Value = -1698.4998
$ws.Range("N134").Value = -21026.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3002
$ws.Range("I63").Value = 1006
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3018
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -2269
$ws.Range("N63").Value = -13498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3002
$ws.Range("I66").Value = 1006
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 9054
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -5310
$ws.Range("N66").Value = -43488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3047.9546
$ws.Range("I80").Value = 2852.2
$ws.Range("J80").Value = 3467.4285
$ws.Range("K80").Value = 2852.2
$ws.Range("L80").Value = 3467.4285
$ws.Range("M80").Value = -1854.2
$ws.Range("N80").Value = -5463.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3047.9546
$ws.Range("I83").Value = 2852.2
$ws.Range("J83").Value = 3467.4285
$ws.Range("K83").Value = 14261
$ws.Range("L83").Value = 17337.1425
$ws.Range("M83").Value = -9269
$ws.Range("N83").Value = -27321.1425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6251975.5
$ws.Range("I122").Value = 10001666
$ws.Range("J122").Value = 2491.5
$ws.Range("K122").Value = 30004998
$ws.Range("L122").Value = 7474.5
$ws.Range("M122").Value = -30002548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 47272.727
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 47272.727
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 47272.727
$ws.Range("N130").Value = -57312.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2174.2942
$ws.Range("I132").Value = 1810.7142
$ws.Range("J132").Value = 3871
$ws.Range("K132").Value = 5432.142599999999
$ws.Range("L132").Value = 11613
$ws.Range("M132").Value = -2902.142599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 19964.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 19964.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 19964.5
$ws.Range("N133").Value = -30084.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2209.0833
$ws.Range("I136").Value = 1958.8422
$ws.Range("J136").Value = 3160
$ws.Range("K136").Value = 5876.5266
$ws.Range("L136").Value = 9480
$ws.Range("M136").Value = -3326.5266
$ws.Range("N136").Value = -14580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1241.625
$ws.Range("I81").Value = 1232.75
$ws.Range("J81").Value = 1250.5
$ws.Range("K81").Value = 2465.5
$ws.Range("L81").Value = 2501
$ws.Range("M81").Value = -1404.5
$ws.Range("N81").Value = -4623

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1241.625
$ws.Range("I84").Value = 1232.75
$ws.Range("J84").Value = 1250.5
$ws.Range("K84").Value = 12327.5
$ws.Range("L84").Value = 12505
$ws.Range("M84").Value = -7023.5
$ws.Range("N84").Value = -23113
